$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price strings so COM does not
# auto-convert them to Number type (source data is stored as text).
$ws.Range("D2").Value = '29.039.28'
$ws.Range("E2").Value = '  -2.09%  '
$ws.Range("D3").Value = '1.908.38'
$ws.Range("E3").Value = '  -4.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.36'
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4604'
$ws.Range("E7").Value = '  -1.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3813'
$ws.Range("E8").Value = '  -3.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07711'
$ws.Range("E9").Value = '  -3.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9781'
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.95'
$ws.Range("E11").Value = '  -4.29%  '
$ws.Range("D12").Value = '1.924.73'
$ws.Range("E12").Value = '  -4.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.949'
$ws.Range("E13").Value = '  -3.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.666'
$ws.Range("E14").Value = '  -3.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07058'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '83.68'
$ws.Range("E17").Value = '  -5.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009487'
$ws.Range("E18").Value = '  -5.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.63'
$ws.Range("E19").Value = '  -4.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").Value = '29.043.80'
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.319'
$ws.Range("E22").Value = '  -4.00%  '
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.095'
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.82'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.11'
$ws.Range("E26").Value = '  -2.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.618'
$ws.Range("E27").Value = '  -3.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.34'
$ws.Range("E28").Value = '  -2.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.849'
$ws.Range("E29").Value = '  -2.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09263'
$ws.Range("E30").Value = '  -2.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8614'
$ws.Range("E31").Value = '  -4.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.088'
$ws.Range("E32").Value = '  -3.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.239'
$ws.Range("E33").Value = '  -7.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.991'
$ws.Range("E34").Value = '  -6.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05682'
$ws.Range("E35").Value = '  -2.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.138'
$ws.Range("E36").Value = '  -3.27%  '
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02035'
$ws.Range("E38").Value = '  -3.39%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5495'
$ws.Range("E39").Value = '  -4.48%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.414'
$ws.Range("E40").Value = '  -6.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1754'
$ws.Range("E41").Value = '  -3.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.314'
$ws.Range("E42").Value = '  -4.77%  '
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5176'
$ws.Range("E44").Value = '  -3.87%  '
$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002681'
$ws.Range("E45").Value = '  -17.63%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.25'
$ws.Range("E46").Value = '  -5.84%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.097'
$ws.Range("E47").Value = '  -3.95%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06796'
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.777'
$ws.Range("E49").Value = '  -3.71%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '109.99'
$ws.Range("E50").Value = '  -3.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  -0.41%  '
